# Scheduled-runner style refresh of per-leve crafting profit figures.
# For each (sheet, row) below, columns H..N hold recomputed market-derived
# numbers (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ]).
# A handful of rows drop/gain a trailing cell entirely (ClearContents mimics
# the sparse-cell removal seen in the source data; plain assignment mimics a
# newly-populated cell) rather than just changing a number in place.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 123.166664
$ws.Range("I8").Value = 107.25
$ws.Range("J8").Value = 155
$ws.Range("K8").Value = 321.75
$ws.Range("L8").Value = 465
$ws.Range("M8").Value = -182.75
$ws.Range("N8").Value = -743

$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents() | Out-Null

$ws.Range("H18").Value = 2639.2
$ws.Range("I18").Value = 2639.2
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 2639.2
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -2355.2

$ws.Range("H19").Value = 924.5
$ws.Range("I19").Value = 966.3333
$ws.Range("J19").Value = 882.6667
$ws.Range("K19").Value = 966.3333
$ws.Range("L19").Value = 882.6667
$ws.Range("M19").Value = -791.3333
$ws.Range("N19").Value = -1232.6667

$ws.Range("H33").Value = 300.83334
$ws.Range("I33").Value = 236.33333
$ws.Range("J33").Value = 365.33334
$ws.Range("K33").Value = 236.33333
$ws.Range("L33").Value = 365.33334
$ws.Range("M33").Value = -7.333329999999989
$ws.Range("N33").Value = -823.33334

$ws.Range("H51").Value = 5491
$ws.Range("I51").Value = 8000.5
$ws.Range("J51").Value = 4933.3335
$ws.Range("K51").Value = 8000.5
$ws.Range("L51").Value = 4933.3335
$ws.Range("M51").Value = -7516.5
$ws.Range("N51").Value = -5901.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 370.27274
$ws.Range("I5").Value = 319.8
$ws.Range("J5").Value = 875
$ws.Range("K5").Value = 319.8
$ws.Range("L5").Value = 875
$ws.Range("M5").Value = -207.8
$ws.Range("N5").Value = -1099

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 370.27274
$ws.Range("I4").Value = 319.8
$ws.Range("J4").Value = 875
$ws.Range("K4").Value = 319.8
$ws.Range("L4").Value = 875
$ws.Range("M4").Value = -204.8
$ws.Range("N4").Value = -1105

$ws.Range("H5").Value = 338.75
$ws.Range("I5").Value = 315.7143
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 315.7143
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = -202.7143
$ws.Range("N5").Value = -726

$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents() | Out-Null

$ws.Range("H38").Value = 19995
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 19995
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 19995
$ws.Range("N38").Value = -20827

$ws.Range("H80").Value = 623.45
$ws.Range("I80").Value = 998.5
$ws.Range("J80").Value = 248.4
$ws.Range("K80").Value = 998.5
$ws.Range("L80").Value = 248.4
$ws.Range("M80").Value = -0.5
$ws.Range("N80").Value = -2244.4

$ws.Range("H83").Value = 623.45
$ws.Range("I83").Value = 998.5
$ws.Range("J83").Value = 248.4
$ws.Range("K83").Value = 4992.5
$ws.Range("L83").Value = 1242
$ws.Range("M83").Value = -0.5
$ws.Range("N83").Value = -11226

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2407.182
$ws.Range("I16").Value = 2166.3333
$ws.Range("J16").Value = 2696.2
$ws.Range("K16").Value = 2166.3333
$ws.Range("L16").Value = 2696.2
$ws.Range("M16").Value = -1879.3333
$ws.Range("N16").Value = -3270.2

$ws.Range("H35").Value = 1171.1428
$ws.Range("I35").Value = 1159.6
$ws.Range("J35").Value = 1200
$ws.Range("K35").Value = 1159.6
$ws.Range("L35").Value = 1200
$ws.Range("M35").Value = -865.5999999999999
$ws.Range("N35").Value = -1788

$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents() | Out-Null

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents() | Out-Null

$ws.Range("H86").Value = 9975.625
$ws.Range("I86").Value = 9959.799999999999
$ws.Range("J86").Value = 10002
$ws.Range("K86").Value = 9959.799999999999
$ws.Range("L86").Value = 10002
$ws.Range("M86").Value = -8836.799999999999
$ws.Range("N86").Value = -12248

$ws.Range("H89").Value = 9975.625
$ws.Range("I89").Value = 9959.799999999999
$ws.Range("J89").Value = 10002
$ws.Range("K89").Value = 49799
$ws.Range("L89").Value = 50010
$ws.Range("M89").Value = -44183
$ws.Range("N89").Value = -61242

$ws.Range("H99").Value = 2673.6
$ws.Range("I99").Value = 1223
$ws.Range("J99").Value = 3640.6667
$ws.Range("K99").Value = 1223
$ws.Range("L99").Value = 3640.6667
$ws.Range("M99").Value = 275
$ws.Range("N99").Value = -6636.6667

$ws.Range("H105").Value = 2933.818
$ws.Range("I105").Value = 1132.75
$ws.Range("J105").Value = 3963
$ws.Range("K105").Value = 1132.75
$ws.Range("L105").Value = 3963
$ws.Range("M105").Value = 614.25
$ws.Range("N105").Value = -7457

$ws.Range("H113").Value = 2407.182
$ws.Range("I113").Value = 2166.3333
$ws.Range("J113").Value = 2696.2
$ws.Range("K113").Value = 2166.3333
$ws.Range("L113").Value = 2696.2
$ws.Range("M113").Value = 3.666700000000219
$ws.Range("N113").Value = -7036.2

$ws.Range("H126").Value = 2673.6
$ws.Range("I126").Value = 1223
$ws.Range("J126").Value = 3640.6667
$ws.Range("K126").Value = 3669
$ws.Range("L126").Value = 10922.0001
$ws.Range("M126").Value = -1199
$ws.Range("N126").Value = -15862.0001

$ws.Range("H132").Value = 1561.625
$ws.Range("I132").Value = 1253.8
$ws.Range("J132").Value = 2661
$ws.Range("K132").Value = 3761.4
$ws.Range("L132").Value = 7983
$ws.Range("M132").Value = -1231.4
$ws.Range("N132").Value = -13043

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 11750
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 11750
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 35250
$ws.Range("N74").Value = -37372

$ws.Range("H77").Value = 11750
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 11750
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 105750
$ws.Range("N77").Value = -116358

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 5000
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 5000
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 5000
$ws.Range("N53").Value = -6262

$ws.Range("H107").Value = 5749.75
$ws.Range("I107").Value = 1500
$ws.Range("J107").Value = 7166.3335
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 7166.3335
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -11006.3335

$ws.Range("H120").Value = 51499.5
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 51499.5
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 51499.5
$ws.Range("N120").Value = -61175.5

$ws.Range("H122").Value = 5257.875
$ws.Range("I122").Value = 4263.1
$ws.Range("J122").Value = 6915.8335
$ws.Range("K122").Value = 12789.3
$ws.Range("L122").Value = 20747.5005
$ws.Range("M122").Value = -10339.3
$ws.Range("N122").Value = -25647.5005

$ws.Range("H126").Value = 3367.7778
$ws.Range("I126").Value = 1786.5
$ws.Range("J126").Value = 4632.8
$ws.Range("K126").Value = 5359.5
$ws.Range("L126").Value = 13898.4
$ws.Range("M126").Value = -2889.5
$ws.Range("N126").Value = -18838.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 126912.5
$ws.Range("I46").Value = 250825
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 250825
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -250637
$ws.Range("N46").Value = -3376

$ws.Range("H132").Value = 2988.5334
$ws.Range("I132").Value = 2222.375
$ws.Range("J132").Value = 3864.1428
$ws.Range("K132").Value = 6667.125
$ws.Range("L132").Value = 11592.4284
$ws.Range("M132").Value = -4137.125
$ws.Range("N132").Value = -16652.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 1205
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1205
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1205
$ws.Range("M17").ClearContents() | Out-Null
$ws.Range("N17").Value = -1549

$ws.Range("H107").Value = 679.7778
$ws.Range("I107").Value = 82.333336
$ws.Range("J107").Value = 978.5
$ws.Range("K107").Value = 247.000008
$ws.Range("L107").Value = 2935.5
$ws.Range("M107").Value = 1672.999992
$ws.Range("N107").Value = -6775.5

$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents() | Out-Null

$ws.Range("H132").Value = 1229.4286
$ws.Range("I132").Value = 1292.2727
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 3876.8181
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -1346.8181
$ws.Range("N132").Value = -8057
